$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.003619074821472
$ws.Range("B1").Value = 1.643280148506165
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.567516803741455
$ws.Range("E1").Value = 1.301176786422729
